$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C (shifts old C..F -> D..G)
$ws.Columns("C:C").Insert()

# Insert a new row before row 4 (shifts old row 4 -> row 5)
$ws.Rows("4:4").Insert()

function Set-Text($addr, $text) {
    # Leading apostrophe forces Excel to store the value as literal text,
    # even when it looks like a number/date/boolean - matches the template's
    # existing "number stored as text" cells (A, C, D, E, F, G columns).
    $ws.Range($addr).Value = "'" + $text
}

# --- Header row ---
Set-Text "C1" "סוג מסמך"
Set-Text "D1" "תיאור"
Set-Text "E1" "מספר מסמך"

# --- Row 2 (ייעוץ ללקוח א) ---
Set-Text "C2" "tax_invoice"

# --- Row 3 (פרויקט עיצוב - updated description + new doc type) ---
Set-Text "C3" "invoice"
Set-Text "D3" "פרויקט עיצוב - חשבונית עסקה"

# --- Row 4 (brand new row) ---
Set-Text "A4" "2024-01-22"
$ws.Range("B4").Value = 11800
Set-Text "C4" "tax_invoice_receipt"
Set-Text "D4" "פרויקט עיצוב - חשבונית מס קבלה"
Set-Text "E4" "1003"
Set-Text "F4" "false"
Set-Text "G4" "paid"

# --- Row 5 (was old row 4 - description/amounts fully replaced) ---
Set-Text "C5" "receipt"
Set-Text "D5" "קבלה על תשלום"
Set-Text "E5" "1004"

# NOTE: the worksheet's <ignoredErrors> "number stored as text" hint range
# (cosmetic-only; no data impact) is preserved as-is by this runtime and is
# not mutable through the exposed Range.Errors COM surface. The sheet's
# <dimension> (A1:G5) does update automatically from the edits above.
